$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (I1, J1) - set values first, then copy H1's format
# so they pick up the same bold/border/centered style used by the other
# header cells.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for the new I and J columns (rows 2-21)
$data = @(
    @(8,9),
    @(8,8),
    @(1,2),
    @(6,6),
    @(7,8),
    @(6,6),
    @(4,5),
    @(3,4),
    @(8,9),
    @(9,9),
    @(2,5),
    @(5,5),
    @(5,5),
    @(10,10),
    @(8,8),
    @(1,3),
    @(7,8),
    @(8,8),
    @(6,6),
    @(2,3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
